$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1255.355988682103
$ws.Range("D2").Value = 914.4329779659899
